$d = $word.ActiveDocument

# 1. Merge the split date text ("...de 20" + "21") into a single run's
#    text "...de 2021". Find/Replace matches across the run boundary and
#    Word collapses the matched range into a single run, carrying over the
#    formatting (Arial / color 000000) of the matched text.
$d.Content.Find.Execute(
    "Guayaquil, ………… de……………………de 2021",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Guayaquil, ………… de……………………de 2021",
    2) | Out-Null

# 2. Remove the now-empty trailing paragraph that used to sit between the
#    date paragraph and the section properties. Since it is the very last
#    paragraph in the body, its mark can't be deleted on its own (Word
#    never leaves the body without a final paragraph mark); instead the
#    range is extended back one character to also swallow the *previous*
#    paragraph's mark, which merges the two paragraphs and removes the
#    empty one.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
if ($last.Range.Text -eq [char]13) {
    $d.Range($last.Range.Start - 1, $last.Range.End).Delete()
}
